$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.085.37"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "'2.455.16"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'509.12"
$ws.Range("E5").Value = "  -2.27%  "
$ws.Range("D6").Value = "'133.93"
$ws.Range("E6").Value = "  +4.38%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.559"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "'2.457.49"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("E13").Value = "  -5.65%  "
$ws.Range("D14").Value = "'2.893.69"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "'57.977.31"
$ws.Range("D16").Value = "'21.96"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("E17").Value = "  +3.11%  "
$ws.Range("D18").Value = "'2.436.54"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "'10.36"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").Value = "'315.50"
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("D22").Value = "'6.45"
$ws.Range("E22").Value = "  +6.26%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "'5.76"
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("D25").Value = "'65.56"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("E28").Value = "  -4.56%  "
$ws.Range("E29").Value = "  +5.44%  "
$ws.Range("D30").Value = "'171.72"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").Value = "'0.0₃0737"
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").Value = "'6.14"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "'1.14"
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").Value = "'18.13"
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("E38").Value = "  +5.14%  "
$ws.Range("D39").Value = "'3.89"
$ws.Range("E39").Value = "  +3.90%  "
$ws.Range("D40").Value = "'36.84"
$ws.Range("E40").Value = "  +1.56%  "
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("D43").Value = "'136.35"
$ws.Range("E43").Value = "  +14.08%  "
$ws.Range("D44").Value = "'3.40"
$ws.Range("D45").Value = "'4.93"
$ws.Range("E45").Value = "  +4.02%  "
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").Value = "'256.44"
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("D48").Value = "'0.0918"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").Value = "'0.0494"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("D51").Value = "'17.28"
$ws.Range("E51").Value = "  +2.13%  "
